$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.00000000000047
$ws.Range("L2").Value = 48.85792322520842
$ws.Range("M2").Value = "[41.1578874972855, 56.55795895313135]"
$ws.Range("P2").Value = 1.62897396852804
$ws.Range("Q2").Value = "[1.452868674633116, 1.8050792624229643]"
$ws.Range("T2").Value = 54.45982482164075
$ws.Range("U2").Value = "[49.43857300476369, 59.48107663851781]"
$ws.Range("X2").Value = 18.51851851851887
$ws.Range("Y2").Value = 17.81781781781815
$ws.Range("Z2").Value = 19.21921921921959

# Row 3
$ws.Range("F3").Value = 25.00000000000047
$ws.Range("H3").Value = [double]"2.752797989558076e-12"
$ws.Range("I3").Value = [double]"2.752797989558076e-12"
$ws.Range("L3").Value = 45.19529806657774
$ws.Range("M3").Value = "[33.50055776892603, 56.890038364229454]"
$ws.Range("N3").Value = [double]"7.130172008373847e-10"
$ws.Range("O3").Value = [double]"7.130172008373847e-10"
$ws.Range("P3").Value = 1.150973885098963
$ws.Range("Q3").Value = "[0.8742369946926551, 1.4277107755052718]"
$ws.Range("R3").Value = [double]"9.812817225451909e-11"
$ws.Range("S3").Value = [double]"9.812817225451909e-11"
$ws.Range("T3").Value = 54.14353255488162
$ws.Range("U3").Value = "[47.71747826320696, 60.56958684655629]"
$ws.Range("X3").Value = 20.4204204204208
$ws.Range("Y3").Value = 19.31931931931968
$ws.Range("Z3").Value = 21.52152152152193

# Row 4
$ws.Range("F4").Value = 25.00000000000047
$ws.Range("H4").Value = [double]"1.110223024625157e-16"
$ws.Range("I4").Value = [double]"1.110223024625157e-16"
$ws.Range("L4").Value = 48.32771471740504
$ws.Range("M4").Value = "[40.789952389444004, 55.86547704536607]"
$ws.Range("N4").Value = [double]"2.220446049250313e-16"
$ws.Range("O4").Value = [double]"2.220446049250313e-16"
$ws.Range("P4").Value = 1.050342288587578
$ws.Range("Q4").Value = "[0.8868159442565782, 1.2138686329185786]"
$ws.Range("R4").Value = [double]"2.220446049250313e-16"
$ws.Range("S4").Value = [double]"2.220446049250313e-16"
$ws.Range("T4").Value = 52.67991592076028
$ws.Range("U4").Value = "[48.6303030722045, 56.72952876931606]"
$ws.Range("X4").Value = 20.82082082082121
$ws.Range("Y4").Value = 20.17017017017055
$ws.Range("Z4").Value = 21.47147147147187

# Row 5
$ws.Range("F5").Value = 25.00000000000047
$ws.Range("H5").Value = [double]"5.551115123125783e-16"
$ws.Range("I5").Value = [double]"5.551115123125783e-16"
$ws.Range("L5").Value = 43.72225516353888
$ws.Range("M5").Value = "[34.97418803100031, 52.470322296077455]"
$ws.Range("N5").Value = [double]"4.225508831723346e-13"
$ws.Range("O5").Value = [double]"4.225508831723346e-13"
$ws.Range("P5").Value = 0.8113422468730391
$ws.Range("Q5").Value = "[0.5975001042863459, 1.0251843894597323]"
$ws.Range("R5").Value = [double]"1.150910478031619e-09"
$ws.Range("S5").Value = [double]"1.150910478031619e-09"
$ws.Range("T5").Value = 54.99131672118721
$ws.Range("U5").Value = "[50.244376856108836, 59.738256586265585]"
$ws.Range("X5").Value = 21.77177177177218
$ws.Range("Y5").Value = 20.92092092092131
$ws.Range("Z5").Value = 22.62262262262305

# Row 6
$ws.Range("F6").Value = 24.01000000000031
$ws.Range("H6").Value = [double]"2.220446049250313e-16"
$ws.Range("I6").Value = [double]"2.220446049250313e-16"
$ws.Range("L6").Value = 47.51840989894028
$ws.Range("M6").Value = "[38.62651464031017, 56.410305157570384]"
$ws.Range("N6").Value = [double]"4.929390229335695e-14"
$ws.Range("O6").Value = [double]"4.929390229335695e-14"
$ws.Range("Q6").Value = "[0.15723686954903737, 0.5597632555945777]"
$ws.Range("R6").Value = 0.000819071048313802
$ws.Range("S6").Value = 0.000819071048313802
$ws.Range("T6").Value = 51.45638648560215
$ws.Range("U6").Value = "[46.38242115799174, 56.53035181321256]"
$ws.Range("X6").Value = 22.64006006006036
$ws.Range("Y6").Value = 21.87097097097126
$ws.Range("Z6").Value = 23.40914914914946

# Row 7
$ws.Range("F7").Value = 24.01000000000031
$ws.Range("H7").Value = [double]"1.110223024625157e-16"
$ws.Range("I7").Value = [double]"1.110223024625157e-16"
$ws.Range("L7").Value = 49.19310030800563
$ws.Range("M7").Value = "[41.05031370718243, 57.33588690882883]"
$ws.Range("N7").Value = [double]"6.661338147750939e-16"
$ws.Range("O7").Value = [double]"6.661338147750939e-16"
$ws.Range("P7").Value = 0.3333421634439633
$ws.Range("Q7").Value = "[0.15723686954904093, 0.5094474573388856]"
$ws.Range("R7").Value = 0.0004159194048691806
$ws.Range("S7").Value = 0.0004159194048691806
$ws.Range("T7").Value = 50.33229824888821
$ws.Range("U7").Value = "[45.753201108015794, 54.911395389760635]"
$ws.Range("X7").Value = 22.73619619619649
$ws.Range("Y7").Value = 22.06324324324353
$ws.Range("Z7").Value = 23.40914914914945

# Row 8
$ws.Range("F8").Value = 24.01000000000031
$ws.Range("L8").Value = 49.83831996008835
$ws.Range("M8").Value = "[40.24957312420541, 59.427066795971285]"
$ws.Range("N8").Value = [double]"1.214583988939921e-13"
$ws.Range("O8").Value = [double]"1.214583988939921e-13"
$ws.Range("P8").Value = 0.3962369112635775
$ws.Range("Q8").Value = "[0.19497371824080734, 0.5975001042863477]"
$ws.Range("R8").Value = 0.000259846526299512
$ws.Range("S8").Value = 0.000259846526299512
$ws.Range("T8").Value = 47.64468133995772
$ws.Range("U8").Value = "[42.46320060572112, 52.82616207419432]"
$ws.Range("X8").Value = 22.49585585585615
$ws.Range("Y8").Value = 21.72676676676705
$ws.Range("Z8").Value = 23.26494494494525

# Row 9
$ws.Range("F9").Value = 24.01000000000031
$ws.Range("H9").Value = [double]"2.164934898019055e-14"
$ws.Range("I9").Value = [double]"2.164934898019055e-14"
$ws.Range("L9").Value = 47.8389806672135
$ws.Range("M9").Value = "[37.56963014763004, 58.10833118679695]"
$ws.Range("N9").Value = [double]"3.683275906496419e-12"
$ws.Range("O9").Value = [double]"3.683275906496419e-12"
$ws.Range("P9").Value = 0.3207632138800394
$ws.Range("Q9").Value = "[0.081763172165501, 0.5597632555945777]"
$ws.Range("R9").Value = 0.009656747921369391
$ws.Range("S9").Value = 0.009656747921369391
$ws.Range("T9").Value = 48.76932268558653
$ws.Range("U9").Value = "[42.96019075048685, 54.57845462068621]"
$ws.Range("X9").Value = 22.78426426426456
$ws.Range("Y9").Value = 21.87097097097126
$ws.Range("Z9").Value = 23.69755755755786
